# Apply "last changes to v1.8.2" edits to the StructureDefinition-TiposVacunaRNI workbook.
#
# 1. Metadata sheet: bump Version 1.8.1 -> 1.8.2 and update the Date timestamp.
# 2. Elements sheet: the root "Extension" row (row 1) gains the same
#    invariant text (ele-1 / ext-1) that already appears on the
#    "Element.extension" row (row 3), in column AJ.

$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item("Metadata")
$metadata.Range("B3").Value = "1.8.2"
$metadata.Range("B8").Value = "2023-09-01T14:45:29-04:00"

$elements = $wb.Worksheets.Item("Elements")
$invariantText = $elements.Range("AJ3").Value2
$elements.Range("AJ1").Value = $invariantText
